# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows for "Betarraga" (Vega Central Mapocho de Santiago)
# at rows 488-491, pushing the existing historical rows down by 4 (488->492 ... 597->601).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (old rows 488:597) down by 4 rows, creating 4 blank rows at 488:491.
$ws.Rows("488:491").Insert()

# New observations to drop into the freshly inserted rows (same template columns
# A,B,C,E,F,G,H,N,Q,R as the rest of the "Betarraga" data block).
$row488 = @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44785, 13, 100114014, "Betarraga", "Sin especificar", "Primera", 3800, 180, 200, 191, "`$/unidad", "Provincia de Melipilla", 191, 1, "Hortaliza")
$row489 = @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44785, 13, 100114014, "Betarraga", "Sin especificar", "Primera", 4300, 170, 180, 175, "`$/unidad", "Región Metropolitana", 175, 1, "Hortaliza")
$row490 = @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44785, 13, 100114014, "Betarraga", "Sin especificar", "Segunda", 5000, 140, 150, 144, "`$/unidad", "Provincia de Melipilla", 144, 1, "Hortaliza")
$row491 = @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44785, 13, 100114014, "Betarraga", "Sin especificar", "Segunda", 9700, 140, 150, 145, "`$/unidad", "Región Metropolitana", 145, 1, "Hortaliza")

$newRows = @($row488, $row489, $row490, $row491)

for ($i = 0; $i -lt 4; $i++) {
    $targetRow = 488 + $i
    $data = $newRows[$i]
    for ($c = 1; $c -le 18; $c++) {
        $ws.Cells.Item($targetRow, $c).Value = $data[$c - 1]
    }
}
